# Update the "kansensya_pcr" style tracking workbook with a new day's
# data (2020-04-29, serial 43950), pushing the footnote rows down by one
# on each of the three data sheets, and bump the "調査中" count text.
#
# NOTE: this COM host requires reads of COM properties to be invoked with
# explicit parens (e.g. `.Value()`), otherwise PowerShell returns the
# property descriptor instead of calling the getter. We avoid relying on
# that here by writing known literal values directly instead of reading
# existing cells back.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "other": insert new data row at 52 (was the footnote row),
# footnote pushed down to row 53, new blank row 54 appended.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

# Row 54 is a brand-new blank row; give it the same formatting as row 53
# (all plain/empty cells) before anything else moves.
$wsOther.Range("A53:I53").Copy()
$wsOther.Range("A54:I54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Push the footnote text from B52 down to B53.
$wsOther.Range("B53").Value = "※他自治体において、3月10日以前の感染者の発生はございません。"

# Copy formatting from the prior data row down onto row 52, then
# overwrite with this day's figures.
$wsOther.Range("A51:I51").Copy()
$wsOther.Range("A52:I52").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsOther.Range("A52").Value = 43950
$wsOther.Range("B52").Value = 0
$wsOther.Range("C52").Value = 11
$wsOther.Range("D52").Value = 4
$wsOther.Range("E52").Value = 3
$wsOther.Range("F52").Value = 1
$wsOther.Range("G52").Value = 0
$wsOther.Range("H52").Value = 7

$wsOther.Range("A52").Select()

# ---------------------------------------------------------------------
# Sheet "kobe": update existing row 76, insert new data row at 77
# (was the footnote row), footnote pushed down to row 78.
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("D76").Value = 3
$wsKobe.Range("E76").Value = 256

# Push the footnote text from B77 down to B78 (A78 stays blank, as A77
# was) - copy A77:B77's current formatting onto the new row 78 first.
$wsKobe.Range("A77:B77").Copy()
$wsKobe.Range("A78:B78").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsKobe.Range("B78").Value = "※24・34・53・58・59・60・158・161・163・192例目は市外在住者です。"

# Copy formatting from the prior data row down onto row 77, then
# overwrite with this day's figures.
$wsKobe.Range("A76:J76").Copy()
$wsKobe.Range("A77:J77").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsKobe.Range("A77").Value = 43950
$wsKobe.Range("B77").Value = 0
$wsKobe.Range("C77").Value = 1854
$wsKobe.Range("D77").Value = 0
$wsKobe.Range("E77").Value = 256
$wsKobe.Range("F77").Value = 126
$wsKobe.Range("G77").Value = 118
$wsKobe.Range("H77").Value = 8
$wsKobe.Range("I77").Value = 4
$wsKobe.Range("J77").Value = 101

$wsKobe.Range("A77").Select()

# ---------------------------------------------------------------------
# Sheet "all": new data row inserted at row 22 (was the first footnote
# row), footnote rows pushed down to 23/24. Processed last so it is
# left as the active sheet/selection on save, matching the source file.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()

# Push the existing footnote text down: old B22 -> B23, old B23 -> B24.
# B24 is a brand-new cell, so copy B23's formatting onto it first (a
# fresh cell would otherwise inherit the plain column default style).
$wsAll.Range("B23").Copy()
$wsAll.Range("B24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsAll.Range("B24").Value = "※15件調査中"
$wsAll.Range("B23").Value = "※24・34・53・58・59・60・158・161・163・192例目は市外在住者です。"

# Copy the formatting of the last real data row down onto row 22,
# then overwrite it with this day's figures.
$wsAll.Range("A21:H21").Copy()
$wsAll.Range("A22:H22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsAll.Range("A22").Value = 43950
$wsAll.Range("B22").Value = 256
$wsAll.Range("C22").Value = 242
$wsAll.Range("D22").Value = 130
$wsAll.Range("E22").Value = 121
$wsAll.Range("F22").Value = 9
$wsAll.Range("G22").Value = 4
$wsAll.Range("H22").Value = 108

$wsAll.Range("B25").Select()
